# Applies the cryptos list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force the cell to stay a text value (Excel would otherwise
    # auto-coerce numeric-looking strings like "111.64" into numbers),
    # then snap the style back to Normal so we don't leave a stray
    # quote-prefixed number format behind on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "43.317.10"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.275.74"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.51%  "
Set-TextCell $ws.Range("D5") "111.64"
$ws.Range("E5").Value = "  +0.98%  "
Set-TextCell $ws.Range("D6") "263.99"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("E8").Value = "  -0.03%  "
Set-TextCell $ws.Range("D9") "0.608"
$ws.Range("E9").Value = "  -0.83%  "
Set-TextCell $ws.Range("D10") "46.57"
$ws.Range("E10").Value = "  -1.80%  "
Set-TextCell $ws.Range("D11") "0.0935"
$ws.Range("E11").Value = "  -0.95%  "
Set-TextCell $ws.Range("D12") "9.25"
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("E13").Value = "  +1.52%  "
Set-TextCell $ws.Range("D14") "15.31"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "2.618.65"
$ws.Range("E15").Value = "  -0.37%  "
Set-TextCell $ws.Range("D16") "0.859"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "2.277.41"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "43.121.31"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -1.02%  "
Set-TextCell $ws.Range("D20") "6.72"
$ws.Range("E20").Value = "  +0.12%  "
Set-TextCell $ws.Range("D21") "72.14"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -1.02%  "
Set-TextCell $ws.Range("D23") "234.03"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("E24").Value = "  +3.44%  "
Set-TextCell $ws.Range("D25") "9.37"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  +1.90%  "
Set-TextCell $ws.Range("D27") "11.31"
$ws.Range("E27").Value = "  -2.22%  "
Set-TextCell $ws.Range("D28") "41.09"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("E29").Value = "  -1.19%  "
Set-TextCell $ws.Range("D30") "2.24"
$ws.Range("E30").Value = "  -0.68%  "
Set-TextCell $ws.Range("D31") "173.50"
$ws.Range("E31").Value = "  -1.40%  "
Set-TextCell $ws.Range("D32") "21.45"
$ws.Range("E32").Value = "  -0.28%  "
Set-TextCell $ws.Range("D33") "0.0894"
$ws.Range("E33").Value = "  -3.18%  "
Set-TextCell $ws.Range("D34") "5.64"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("E35").Value = "  +3.60%  "
Set-TextCell $ws.Range("D36") "0.0380"
$ws.Range("E36").Value = "  +5.56%  "
Set-TextCell $ws.Range("D37") "4.67"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  -3.42%  "
Set-TextCell $ws.Range("D40") "2.57"
$ws.Range("E40").Value = "  +6.87%  "
Set-TextCell $ws.Range("D41") "14.22"
$ws.Range("E41").Value = "  +4.81%  "
Set-TextCell $ws.Range("D42") "74.83"
$ws.Range("E42").Value = "  +3.71%  "
Set-TextCell $ws.Range("D43") "0.235"
$ws.Range("E43").Value = "  -2.73%  "
Set-TextCell $ws.Range("D44") "6.07"
$ws.Range("E44").Value = "  -2.85%  "
Set-TextCell $ws.Range("D45") "0.999"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("E47").Value = "  +4.89%  "
Set-TextCell $ws.Range("D48") "8.54"
$ws.Range("E48").Value = "  -2.96%  "
Set-TextCell $ws.Range("D49") "0.0989"
$ws.Range("E49").Value = "  -0.93%  "
Set-TextCell $ws.Range("D50") "100.27"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws.Range("D51") "0.597"
$ws.Range("E51").Value = "  +9.93%  "
